# Finished Week 13 logging
$wb = $excel.ActiveWorkbook

# OFF sheet - Row 3 (Road) updated target depth totals
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 172
$wsOff.Range("C3").Value = 125
$wsOff.Range("D3").Value = 39
$wsOff.Range("E3").Value = 17
$wsOff.Range("F3").Value = 4

# DEF sheet - Row 3 (Road) updated target depth totals
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 155
$wsDef.Range("C3").Value = 92
$wsDef.Range("D3").Value = 45
$wsDef.Range("E3").Value = 22
$wsDef.Range("F3").Value = 3
